$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Aline Lima"
$ws.Range("D4").Value = "Jane"
$ws.Range("F4").Value = "Graca"
$ws.Range("G4").Value = "Lurdes"
$ws.Range("J4").Value = "Vinicius"
$ws.Range("K4").Value = "Antonio"
$ws.Range("D5").Value = "Keila"
$ws.Range("G5").Value = "Patricia Rodrigues"
$ws.Range("C6").Value = "Ediane"
$ws.Range("D6").Value = "Lucia"
$ws.Range("F6").Value = "Conceicao"
$ws.Range("G6").Value = "Patricia Dias"
$ws.Range("J6").Value = "Rodolfo"
$ws.Range("K6").Value = "Antonio"
$ws.Range("C8").Value = "Edith"
$ws.Range("D8").Value = "Isabele"
$ws.Range("G8").Value = "Valquiria"
$ws.Range("K8").Value = "Antonio"
$ws.Range("M8").Value = "Clayton"
$ws.Range("K9").Value = "Telma"
$ws.Range("M9").Value = "Douglas Oliveira"
$ws.Range("C11").Value = "Lucia"
$ws.Range("D11").Value = "Helaine Camilo"
$ws.Range("F11").Value = "Alana"
$ws.Range("J11").Value = "Keila"
$ws.Range("K11").Value = "Valquiria"
$ws.Range("L11").Value = "Marcio"
$ws.Range("D12").Value = "Vanda"
$ws.Range("G12").Value = "Alana"
$ws.Range("K12").Value = "Rodolfo"
$ws.Range("C13").Value = "Rodolfo"
$ws.Range("D13").Value = "Robson"
$ws.Range("F13").Value = "Lindoia"
$ws.Range("G13").Value = "Valquiria"
$ws.Range("K13").Value = "Lurdes"
$ws.Range("D14").Value = "Lucia"
$ws.Range("C15").Value = "Daiana"
$ws.Range("D15").Value = "Eliane"
$ws.Range("J15").Value = "Lucia"
$ws.Range("K15").Value = "Keila"
$ws.Range("L15").Value = "Rodolfo"
$ws.Range("K16").Value = "Alex"
$ws.Range("M16").Value = "Eliane"
$ws.Range("C18").Value = "Lurdes"
$ws.Range("F18").Value = "Aline Lima"
$ws.Range("G18").Value = "Graca"
$ws.Range("J18").Value = "Icaro"
$ws.Range("K18").Value = "Antonio"
$ws.Range("L18").Value = "Dario"
$ws.Range("M18").Value = "EMPTY"
$ws.Range("D19").Value = "Patricia Dias"
$ws.Range("G19").Value = "Valquiria"
$ws.Range("J19").Value = "Keila"
$ws.Range("D20").Value = "Helaine Camilo"
$ws.Range("F20").Value = "Edith"
$ws.Range("G20").Value = "Patricia Dias"
$ws.Range("J20").Value = "Antonio"
$ws.Range("K20").Value = "EMPTY"
$ws.Range("J22").Value = "Aline Lima"
$ws.Range("K22").Value = "Beth"
$ws.Range("L22").Value = "Vinicius"
$ws.Range("K23").Value = "Karol"
$ws.Range("M23").Value = "Marcio"
$ws.Range("D25").Value = "Alana"
$ws.Range("F25").Value = "Lucia"
$ws.Range("J25").Value = "Eliane"
$ws.Range("L25").Value = "Amintas"
$ws.Range("M25").Value = "Marcio"
$ws.Range("D26").Value = "Alana"
$ws.Range("G26").Value = "Patricia Rodrigues"
$ws.Range("J26").Value = "Jessica Silva"
$ws.Range("K26").Value = "EMPTY"
$ws.Range("C27").Value = "Rodolfo"
$ws.Range("D27").Value = "Lindoia"
$ws.Range("F27").Value = "Lurdes"
$ws.Range("G27").Value = "Lucia"
$ws.Range("D28").Value = "Vanda"
$ws.Range("C29").Value = "Keila"
$ws.Range("D29").Value = "Daiana"
$ws.Range("G29").Value = "Lucia"
$ws.Range("J29").Value = "Lurdes"
$ws.Range("K29").Value = "EMPTY"
$ws.Range("L29").Value = "Rodolfo"
$ws.Range("M29").Value = "Clayton"
$ws.Range("K30").Value = "Eliane"
$ws.Range("C32").Value = "Daniel"
$ws.Range("D32").Value = "EMPTY"
$ws.Range("F32").Value = "Aline Lima"
$ws.Range("G32").Value = "Lurdes"
$ws.Range("J32").Value = "Patricia Dias"
$ws.Range("K32").Value = "EMPTY"
$ws.Range("L32").Value = "Dario"
$ws.Range("M32").Value = "Icaro"

# December 31st row is removed; this shifts the "Data de geracao" footer
# block (and the blank rows around it) up by one row, and drops the last
# trailing blank row, matching the new dimension A1:M36.
$ws.Rows("33").Delete()

# Update the generation timestamp (now on row 34 after the shift above).
$ws.Range("C34").Value = 43818.59197647273
